$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: place/coordinate pair changes from "unna"/"51.5333,7.6833" to "brugge"/"51.2089,3.2242"
$ws.Range("A2").Value = "brugge"
$ws.Range("B2").Value = "51.2089,3.2242"

# New rows 3-7 with additional place / coordinate pairs
$ws.Range("A3").Value = "gent"
$ws.Range("B3").Value = "51.05,3.7167"

$ws.Range("A4").Value = "dusseldorf"
$ws.Range("B4").Value = "51.2217,6.7762"

$ws.Range("A5").Value = "genk"
$ws.Range("B5").Value = "50.965,5.5008"

$ws.Range("A6").Value = "brussels"
$ws.Range("B6").Value = "50.8504,4.3488"

$ws.Range("A7").Value = "antwerp"
$ws.Range("B7").Value = "51.2199,4.4035"

# Re-apply explicit column widths (keeps custom-width flag true on all 3 columns,
# widens column B to fit the newly added, longer coordinate strings)
$ws.Columns.Item(1).ColumnWidth = 43.59
$ws.Columns.Item(2).ColumnWidth = 14.14
$ws.Columns.Item(3).ColumnWidth = 19.75

# Move selection to the newly active/last-entered cell
$ws.Range("A7").Select()
